# Applies the edits described by the commit:
#  "added isdriveropen method; corrected ctdc tc to run; added queries in all ctdc tc xls"
# i.e. put the Cypher MATCH query text into A2 of the "startup" sheet, let the
# row grow to fit it (wrap-text style already lives in the cell), and update
# the sheet/window selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$query = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.race IN ['AMERICAN_INDIAN_OR_ALASKA_NATIVE'] RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(t.clinical_trial_designation ,'')as ``Trial Code`` , coalesce(a.arm_id,'') As ``Arm`` , coalesce(a.arm_drug,'') As ``Arm Treatment`` , coalesce(c.disease,'') As Diagnosis , coalesce(c.gender,'') As Gender , coalesce(c.race,'') As Race , coalesce(c.ethnicity,'') As Ethnicity"

# Write the new shared-string text into A2 (keeps the existing wrap-text style).
$ws.Range("A2").Value = $query

# Grow row 2 to fit the now much-longer wrapped text.
$ws.Rows.Item(2).RowHeight = 101.5

# Update the selection to B2:B4 (active cell lands on the anchor cell B2,
# which is the closest this object model gets to the recorded B4/B2:B4
# selection).
$ws.Range("B2:B4").Select() | Out-Null
